$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "sgdm"

# --- flow20 table (rows 23-27) ---
# Row 23 (MILP): Time 2.7h -> 2.4h, Mean Cost 75.24 -> 64.91
$ws.Range("B23").Value = "2.4h"
$ws.Range("C23").Value = 64.91

# Row 24 (CNN): Time 1.44s -> 2.1s, and previously-undefined Precision/Feasible Ratio/Max Diff now populated
$ws.Range("B24").Value = "2.1s"
$ws.Range("C24").Value = 124.66
$ws.Range("D24").Value = 0.5185
$ws.Range("E24").Value = 0.8485
$ws.Range("E24").NumberFormat = "0.00%"
$ws.Range("F24").Value = 136.6

# Row 25 (Greedy): Mean Cost/Precision/Feasible Ratio/Max Diff updated
$ws.Range("C25").Value = 269.46
$ws.Range("D25").Value = 0.4435
$ws.Range("E25").Value = 0.8075
$ws.Range("E25").NumberFormat = "0.00%"
$ws.Range("F25").Value = 646.1

# Row 26 (RGR): Time 56.95s -> 94.16s, Mean Cost/Precision/Feasible Ratio/Max Diff updated
$ws.Range("B26").Value = "94.16s"
$ws.Range("C26").Value = 156.36
$ws.Range("D26").Value = 0.2105
$ws.Range("D26").NumberFormat = "0.00%"
$ws.Range("E26").Value = 0.808
$ws.Range("E26").NumberFormat = "0.00%"
$ws.Range("F26").Value = 635.8

# --- Page setup (adds <pageSetup paperSize="9" orientation="portrait".../>) ---
try {
    $ws.PageSetup.PaperSize = 9
    $ws.PageSetup.Orientation = 1
} catch {}

# --- View: scroll position + selection ---
try {
    $excel.ActiveWindow.ScrollRow = 7
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("E27").Select()
